$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Now go into myCourses," - remove the spell-check proofErr markup
#    and re-flow the runs so the paragraph reads the same but without
#    the <w:proofErr/> wrapped "myCourses" split.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Now go into myCourses,") | Out-Null
$startPos = $rng.Start
$endPos = $rng.End
$delRng = $d.Range($startPos, $endPos)
$delRng.Delete()

$insA = $d.Range($startPos, $startPos)
$insA.InsertBefore("yCourses,")

$insB = $d.Range($startPos, $startPos)
$insB.InsertBefore("Now go into m")

# ---------------------------------------------------------------------
# 2) Insert " with your neighbors in class" right before the final
#    period of "... Feel free to make these determinations as a team."
#    and move the singleton _GoBack bookmark to sit right after the
#    newly typed text (its last-edit-location), exactly where Word
#    would leave it after a live edit. A short-lived helper bookmark
#    is used to pin the run boundary between "team" and the new text
#    so the paragraph keeps three separate runs once the edit settles.
# ---------------------------------------------------------------------
$teamRng = $d.Content
$teamRng.Find.Execute("make these determinations as a team.") | Out-Null
$periodPos = $teamRng.End - 1

$newText = " with your neighbors in class"
$newTextRng = $d.Range($periodPos, $periodPos)
$newTextRng.InsertBefore($newText)

# pin the run boundary between "team" and the newly typed text so the
# paragraph keeps the new text as its own run once the edit settles
$splitPos = $d.Range($periodPos, $periodPos)
$d.Bookmarks.Add("ZZZTempSplit", $splitPos)

$goBackAt = $periodPos + $newText.Length
$goBackRng = $d.Range($goBackAt, $goBackAt)
$d.Bookmarks.Add("_GoBack", $goBackRng)

if ($d.Bookmarks.Exists("ZZZTempSplit")) {
    $d.Bookmarks("ZZZTempSplit").Delete()
}

# ---------------------------------------------------------------------
# 3) "...appropriate dropbox before next class." - remove the
#    spell-check proofErr markup around "dropbox".
# ---------------------------------------------------------------------
$d.Content.Find.Execute(", to the appropriate dropbox before next class.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    ", to the appropriate dropbox before next class.", 2) | Out-Null
